$wb = $excel.ActiveWorkbook

# --- Fix typo in "Materias" sheet header: "tirno" -> "turno" ---
$wsMaterias = $wb.Worksheets.Item("Materias")
$wsMaterias.Range("D1").Value = "turno"

# --- Update remembered selection on "Materias" sheet to D8 ---
[void]$wsMaterias.Range("D8").Select()

# --- Make "Grupos" the active/selected sheet (was "Asignaciones") ---
$wsGrupos = $wb.Worksheets.Item("Grupos")
[void]$wsGrupos.Activate()
